# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2135.3333
$ws.Range("J17").Value = 292.46155
$ws.Range("L17").Value = 877.38465
$ws.Range("N17").Value = -1213.38465
$ws.Range("H92").Value = 2000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").Value = 2000
$ws.Range("N92").Value = -4496
$ws.Range("H100").Value = 1462.7778
$ws.Range("I100").Value = 1544.6154
$ws.Range("J100").Value = 1386.7858
$ws.Range("K100").Value = 1544.6154
$ws.Range("L100").Value = 1386.7858
$ws.Range("M100").Value = -1003.6154
$ws.Range("N100").Value = -2468.7858
$ws.Range("H112").Value = 25001176
$ws.Range("I112").Value = 525
$ws.Range("J112").Value = 26317000
$ws.Range("K112").Value = 1575
$ws.Range("L112").Value = 78951000
$ws.Range("M112").Value = -467
$ws.Range("N112").Value = -78953216
$ws.Range("H141").Value = 5454.909
$ws.Range("I141").Value = 5444.3335
$ws.Range("J141").Value = 5502.5
$ws.Range("K141").Value = 16333.0005
$ws.Range("L141").Value = 16507.5
$ws.Range("M141").Value = -11153.0005
$ws.Range("N141").Value = -26867.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 80.42856999999999
$ws.Range("I5").Value = 75.25
$ws.Range("K5").Value = 75.25
$ws.Range("M5").Value = 36.75

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 80.42856999999999
$ws.Range("I4").Value = 75.25
$ws.Range("K4").Value = 75.25
$ws.Range("M4").Value = 39.75
$ws.Range("H134").Value = 43909.395
$ws.Range("I134").Value = 43909.395
$ws.Range("K134").Value = 131728.185
$ws.Range("M134").Value = -129193.185

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1238.2727
$ws.Range("I16").Value = 1112.2
$ws.Range("J16").Value = 2499
$ws.Range("K16").Value = 1112.2
$ws.Range("L16").Value = 2499
$ws.Range("M16").Value = -825.2
$ws.Range("N16").Value = -3073
$ws.Range("H22").Value = 243.875
$ws.Range("I22").Value = 215.5
$ws.Range("K22").Value = 215.5
$ws.Range("M22").Value = 134.5
$ws.Range("H62").Value = 3019
$ws.Range("I62").Value = 3101.2856
$ws.Range("J62").Value = 2875
$ws.Range("K62").Value = 3101.2856
$ws.Range("L62").Value = 2875
$ws.Range("M62").Value = -2477.2856
$ws.Range("N62").Value = -4123
$ws.Range("H65").Value = 3019
$ws.Range("I65").Value = 3101.2856
$ws.Range("J65").Value = 2875
$ws.Range("K65").Value = 15506.428
$ws.Range("L65").Value = 14375
$ws.Range("M65").Value = -12386.428
$ws.Range("N65").Value = -20615
$ws.Range("H99").Value = 1115.875
$ws.Range("I99").Value = 1115.875
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1115.875
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 382.125
$ws.Range("H107").Value = 339.39026
$ws.Range("I107").Value = 351.21875
$ws.Range("J107").Value = 297.33334
$ws.Range("K107").Value = 351.21875
$ws.Range("L107").Value = 297.33334
$ws.Range("M107").Value = 1568.78125
$ws.Range("N107").Value = -4137.33334
$ws.Range("H113").Value = 1238.2727
$ws.Range("I113").Value = 1112.2
$ws.Range("J113").Value = 2499
$ws.Range("K113").Value = 1112.2
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = 1057.8
$ws.Range("N113").Value = -6839
$ws.Range("H122").Value = 9617033
$ws.Range("I122").Value = 17858584
$ws.Range("J122").Value = 1888.8334
$ws.Range("K122").Value = 53575752
$ws.Range("L122").Value = 5666.5002
$ws.Range("M122").Value = -53573302
$ws.Range("N122").Value = -10566.5002
$ws.Range("H126").Value = 1115.875
$ws.Range("I126").Value = 1115.875
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3347.625
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -877.625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 986201.25
$ws.Range("I2").Value = 3.3333333
$ws.Range("J2").Value = 1831513.8
$ws.Range("K2").Value = 19.9999998
$ws.Range("L2").Value = 10989082.8
$ws.Range("M2").Value = 93.0000002
$ws.Range("N2").Value = -10989308.8
$ws.Range("H12").Value = 62.64
$ws.Range("I12").Value = 86.36364
$ws.Range("J12").Value = 44
$ws.Range("K12").Value = 259.09092
$ws.Range("L12").Value = 132
$ws.Range("M12").Value = -86.09091999999998
$ws.Range("N12").Value = -478
$ws.Range("H116").Value = 1095.3334
$ws.Range("I116").Value = 365.6
$ws.Range("J116").Value = 2007.5
$ws.Range("K116").Value = 1096.8
$ws.Range("L116").Value = 6022.5
$ws.Range("M116").Value = 2345.2
$ws.Range("N116").Value = -12906.5
$ws.Range("H131").Value = 1221774.1
$ws.Range("I131").Value = 4894.231
$ws.Range("J131").Value = 1451041.4
$ws.Range("K131").Value = 14682.693
$ws.Range("L131").Value = 4353124.199999999
$ws.Range("M131").Value = -9642.692999999999
$ws.Range("N131").Value = -4363204.199999999
$ws.Range("H133").Value = 6172.577
$ws.Range("I133").Value = 3067.8
$ws.Range("J133").Value = 6911.8096
$ws.Range("K133").Value = 9203.400000000001
$ws.Range("L133").Value = 20735.4288
$ws.Range("M133").Value = -4143.400000000001
$ws.Range("N133").Value = -30855.4288

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 6013.5
$ws.Range("I48").Value = 6013.5
$ws.Range("K48").Value = 6013.5
$ws.Range("M48").Value = -5528.5
$ws.Range("H80").Value = 135212.33
$ws.Range("I80").Value = 2452.5
$ws.Range("J80").Value = 173143.72
$ws.Range("K80").Value = 2452.5
$ws.Range("L80").Value = 173143.72
$ws.Range("M80").Value = -1454.5
$ws.Range("N80").Value = -175139.72
$ws.Range("H83").Value = 135212.33
$ws.Range("I83").Value = 2452.5
$ws.Range("J83").Value = 173143.72
$ws.Range("K83").Value = 12262.5
$ws.Range("L83").Value = 865718.6
$ws.Range("M83").Value = -7270.5
$ws.Range("N83").Value = -875702.6
$ws.Range("H113").Value = 1944.1818
$ws.Range("I113").Value = 2057.3
$ws.Range("J113").Value = 813
$ws.Range("K113").Value = 2057.3
$ws.Range("L113").Value = 813
$ws.Range("M113").Value = 112.6999999999998
$ws.Range("N113").Value = -5153
$ws.Range("H126").Value = 2307.2
$ws.Range("I126").Value = 2687.7144
$ws.Range("J126").Value = 1419.3334
$ws.Range("K126").Value = 8063.1432
$ws.Range("L126").Value = 4258.0002
$ws.Range("M126").Value = -5593.1432
$ws.Range("N126").Value = -9198.0002

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 10966.667
$ws.Range("J59").Value = 10966.667
$ws.Range("L59").Value = 10966.667
$ws.Range("N59").Value = -12274.667
$ws.Range("H74").Value = 33450
$ws.Range("J74").Value = 33450
$ws.Range("L74").Value = 33450
$ws.Range("N74").Value = -35446
$ws.Range("H77").Value = 33450
$ws.Range("J77").Value = 33450
$ws.Range("L77").Value = 100350
$ws.Range("N77").Value = -110334
$ws.Range("H93").Value = 2511.9443
$ws.Range("I93").Value = 1853.6154
$ws.Range("J93").Value = 4223.6
$ws.Range("K93").Value = 1853.6154
$ws.Range("L93").Value = 4223.6
$ws.Range("M93").Value = -605.6153999999999
$ws.Range("N93").Value = -6719.6
$ws.Range("H111").Value = 36374.57
$ws.Range("J111").Value = 36374.57
$ws.Range("L111").Value = 36374.57
$ws.Range("N111").Value = -44554.57

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 12899.333
$ws.Range("J49").Value = 12899.333
$ws.Range("L49").Value = 12899.333
$ws.Range("N49").Value = -13359.333
$ws.Range("H62").Value = 14329.833
$ws.Range("I62").Value = 2980
$ws.Range("J62").Value = 16599.8
$ws.Range("K62").Value = 2980
$ws.Range("L62").Value = 16599.8
$ws.Range("M62").Value = -2356
$ws.Range("N62").Value = -17847.8
$ws.Range("H65").Value = 14329.833
$ws.Range("I65").Value = 2980
$ws.Range("J65").Value = 16599.8
$ws.Range("K65").Value = 14900
$ws.Range("L65").Value = 82999
$ws.Range("M65").Value = -11780
$ws.Range("N65").Value = -89239
